$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "last_name" header (C1) was renamed to "password"
$ws.Range("C1").Value = "password"

# Update the current selection to reflect where the user left the cursor
[void]$ws.Range("H15").Select()
